# Refresh the cryptocurrency price/volume table (and reorder two swapped rows)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.919.78"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.66%  "
$ws.Range("D3").Value = "'3.721.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'610.82"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.94%  "
$ws.Range("D6").Value = "'187.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.38%  "
$ws.Range("D7").Value = "'3.717.71"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.94%  "
$ws.Range("D8").Value = "'0.638"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.86%  "
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("D10").Value = "'0.717"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("D11").Value = "'0.161"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.55%  "
$ws.Range("D12").Value = "'57.09"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +9.15%  "
$ws.Range("D13").Value = "'0.0000290"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.17%  "
$ws.Range("D14").Value = "'10.60"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.07%  "
$ws.Range("D15").Value = "'4.331.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.15%  "
$ws.Range("D16").Value = "'3.742.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("D17").Value = "'19.32"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "'13.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D19").Value = "'0.126"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("B20").Value = "Polygon"
$ws.Range("C20").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D20").Value = "'1.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").Value = "'68.825.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.47%  "
$ws.Range("D22").Value = "'411.65"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.79%  "
$ws.Range("D23").Value = "'4.59"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.68%  "
$ws.Range("D24").Value = "'89.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.82%  "
$ws.Range("D25").Value = "'3.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.77%  "
$ws.Range("D26").Value = "'12.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.56%  "
$ws.Range("D27").Value = "'11.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.72%  "
$ws.Range("D28").Value = "'6.05"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.99%  "
$ws.Range("D29").Value = "'3.79"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.70%  "
$ws.Range("D30").Value = "'9.63"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.93%  "
$ws.Range("D31").Value = "'33.22"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("D32").Value = "'7.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -9.02%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.124"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.38%  "
$ws.Range("B34").Value = "Cosmos"
$ws.Range("C34").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D34").Value = "'12.63"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("D35").Value = "'626.55"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.70%  "
$ws.Range("D36").Value = "'44.67"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.26%  "
$ws.Range("D37").Value = "'65.98"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("D38").Value = "'0.0₃0834"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -9.88%  "
$ws.Range("D39").Value = "'0.415"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.81%  "
$ws.Range("D40").Value = "'0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").Value = "'0.140"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.95%  "
$ws.Range("D43").Value = "'3.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("D44").Value = "'0.0444"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.64%  "
$ws.Range("D45").Value = "'2.63"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.02%  "
$ws.Range("D46").Value = "'0.140"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.69%  "
$ws.Range("D47").Value = "'2.857.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.05%  "
$ws.Range("D48").Value = "'2.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.36%  "
$ws.Range("D49").Value = "'9.14"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.23%  "
$ws.Range("D50").Value = "'2.68"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -19.29%  "
$ws.Range("D51").Value = "'3.11"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.05%  "
